# The workbook tracked "TMT" as the only supported labelling technology;
# this commit generalises the wording/sheets to "Label" (so e.g. droplet-based
# techniques that aren't TMT can use the same mapping sheets), and also
# updates two README guidance cells whose text referred to TMT specifically.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename the two mapping sheets.
#    Note the "Template" sheet's new name purposely has NO space before
#    "Template" ("Well to Label mapping -Template").
# ---------------------------------------------------------------------------
$wsDefault = $wb.Worksheets.Item("Well to TMT mapping - Default")
$wsTemplate = $wb.Worksheets.Item("Well to TMT mapping - Template")

$wsDefault.Name = "Well to Label mapping - Default"
$wsTemplate.Name = "Well to Label mapping -Template"

# ---------------------------------------------------------------------------
# 2) README sheet text updates.
# ---------------------------------------------------------------------------
$wsReadme = $wb.Worksheets.Item("README")

# "Labels file" -> "Labels file / droplet location file"
$wsReadme.Range("B8").Value = "Labels file / droplet location file"

# Missing-cell-data guidance no longer singles out TMT labels
$wsReadme.Range("C18").Value = "Assign a value to wells that are missing cell data"

# Well to TMT mapping CSV guidance -> Well to Label mapping CSV guidance
$wsReadme.Range("B20").Value = 'Well to Label mapping CSV (Only if "Label-based" selected as technology)'
$wsReadme.Range("C20").Value = "Choose the mapping of well to Label, if 'Default' selected 'Well to Label mapping - Default' is used, else 'Well to Label mapping - Template' can be populated and uploaded as a csv"

# ---------------------------------------------------------------------------
# 3) "Well to Label mapping - Default" sheet: header TMT -> Label, and the
#    two previously-blank "Well" cells (for Carrier/Reference reporter rows)
#    now explicitly read "Empty". The sheet also loses its per-cell styling
#    (rows 1-19 drop their style indices back to the default/general style).
# ---------------------------------------------------------------------------
$wsDefault.Range("B1").Value = "Label"
$wsDefault.Range("A18").Value = "Empty"
$wsDefault.Range("A19").Value = "Empty"
$wsDefault.Range("A1:B19").ClearFormats()

# ---------------------------------------------------------------------------
# 4) "Well to Label mapping -Template" sheet: header TMT -> Label.
# ---------------------------------------------------------------------------
$wsTemplate.Range("B1").Value = "Label"
